$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: shift existing data (years 2021..2000) down by one row, columns A:M only.
# Column O (Fonte/source links) must stay fixed in its rows.
$src = $ws.Range("A2:M23").Value2
$ws.Range("A3:M24").Value2 = $src

# Step 2: write the new row 2 for year 2022.
$ws.Range("A2").Value = 2022
$ws.Range("B2").Value = 0.55000000000000004
$ws.Range("C2:M2").Value = ""

Write-Output "done"
